$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.626.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.861.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.37%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4685"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.49%  "

$ws.Range("E8").Value = "  -0.89%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.80"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.38%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07982"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.12%  "

$ws.Range("E11").Value = "  -2.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.76"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.883.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.988"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.243"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.013"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06732"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001043"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.60%  "

$ws.Range("E20").Value = "  -1.10%  "

$ws.Range("E21").Value = "  +0.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.604.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.453"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.64%  "

$ws.Range("E24").Value = "  -1.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.306"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.095.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.138"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.426"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9756"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.78%  "

$ws.Range("E33").Value = "  -0.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.614"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.289"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.334"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06050"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02230"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.284"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.191"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.25%  "

$ws.Range("E41").Value = "  +0.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5936"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1881"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.94%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.251"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5629"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.198"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06762"
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.16%  "
